$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3337.5
$ws.Range("I51").Value = 4950
$ws.Range("J51").Value = 2800
$ws.Range("K51").Value = 4950
$ws.Range("L51").Value = 2800
$ws.Range("M51").Value = -4466
$ws.Range("N51").Value = -3768

$ws.Range("H132").Value = 7938029
$ws.Range("I132").Value = 9525368
$ws.Range("J132").Value = 1332
$ws.Range("K132").Value = 28576104
$ws.Range("L132").Value = 3996
$ws.Range("M132").Value = -28573574
$ws.Range("N132").Value = -9056

$ws.Range("H137").Value = 2009.7858
$ws.Range("I137").Value = 1818.4
$ws.Range("K137").Value = 5455.200000000001
$ws.Range("M137").Value = -2905.200000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1558.5
$ws.Range("I2").Value = 1198.5
$ws.Range("K2").Value = 1198.5
$ws.Range("M2").Value = -1085.5

$ws.Range("H61").Value = 2128.24
$ws.Range("I61").Value = 1466
$ws.Range("J61").Value = 3305.5557
$ws.Range("K61").Value = 1466
$ws.Range("L61").Value = 3305.5557
$ws.Range("M61").Value = -1254
$ws.Range("N61").Value = -3729.5557

$ws.Range("H74").Value = 982.3077
$ws.Range("I74").Value = 931.4
$ws.Range("J74").Value = 1152
$ws.Range("K74").Value = 931.4
$ws.Range("L74").Value = 1152
$ws.Range("M74").Value = -57.39999999999998
$ws.Range("N74").Value = -2900

$ws.Range("H77").Value = 982.3077
$ws.Range("I77").Value = 931.4
$ws.Range("J77").Value = 1152
$ws.Range("K77").Value = 4657
$ws.Range("L77").Value = 5760
$ws.Range("M77").Value = -289
$ws.Range("N77").Value = -14496

$ws.Range("H116").Value = 1558.5
$ws.Range("I116").Value = 1198.5
$ws.Range("K116").Value = 1198.5
$ws.Range("M116").Value = 1095.5

$ws.Range("H132").Value = 7825.222
$ws.Range("I132").Value = 11987.8
$ws.Range("K132").Value = 35963.39999999999
$ws.Range("M132").Value = -33433.39999999999

$ws.Range("H136").Value = 2128.24
$ws.Range("I136").Value = 1466
$ws.Range("J136").Value = 3305.5557
$ws.Range("K136").Value = 4398
$ws.Range("L136").Value = 9916.667099999999
$ws.Range("M136").Value = -1848
$ws.Range("N136").Value = -15016.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1558.5
$ws.Range("I3").Value = 1198.5
$ws.Range("K3").Value = 1198.5
$ws.Range("M3").Value = -1084.5

$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").Value = ""

$ws.Range("H134").Value = 16847.166
$ws.Range("I134").Value = 22331.521
$ws.Range("J134").Value = 2222.2222
$ws.Range("K134").Value = 66994.56299999999
$ws.Range("L134").Value = 6666.6666
$ws.Range("M134").Value = -64459.56299999999
$ws.Range("N134").Value = -11736.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5715812.5
$ws.Range("I31").Value = 1495.6923
$ws.Range("J31").Value = 22223838
$ws.Range("K31").Value = 1495.6923
$ws.Range("L31").Value = 22223838
$ws.Range("M31").Value = -1200.6923
$ws.Range("N31").Value = -22224428

$ws.Range("H34").Value = 5715812.5
$ws.Range("I34").Value = 1495.6923
$ws.Range("J34").Value = 22223838
$ws.Range("K34").Value = 1495.6923
$ws.Range("L34").Value = 22223838
$ws.Range("M34").Value = -1293.6923
$ws.Range("N34").Value = -22224242

$ws.Range("H58").Value = 1097.1428
$ws.Range("I58").Value = 1035.3846
$ws.Range("J58").Value = 1900
$ws.Range("K58").Value = 1035.3846
$ws.Range("L58").Value = 1900
$ws.Range("M58").Value = -832.3846000000001
$ws.Range("N58").Value = -2306

$ws.Range("H132").Value = 2099.5476
$ws.Range("I132").Value = 2004
$ws.Range("J132").Value = 2254.8125
$ws.Range("K132").Value = 6012
$ws.Range("L132").Value = 6764.4375
$ws.Range("M132").Value = -3482
$ws.Range("N132").Value = -11824.4375

$ws.Range("H134").Value = 936.881
$ws.Range("I134").Value = 908.9487
$ws.Range("J134").Value = 1300
$ws.Range("K134").Value = 2726.8461
$ws.Range("L134").Value = 3900
$ws.Range("M134").Value = -191.8461000000002
$ws.Range("N134").Value = -8970

$ws.Range("H136").Value = 1097.1428
$ws.Range("I136").Value = 1035.3846
$ws.Range("J136").Value = 1900
$ws.Range("K136").Value = 3106.1538
$ws.Range("L136").Value = 5700
$ws.Range("M136").Value = -556.1538
$ws.Range("N136").Value = -10800

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3308.8696
$ws.Range("I3").Value = 1760.2667
$ws.Range("K3").Value = 5280.800099999999
$ws.Range("M3").Value = -5168.800099999999

$ws.Range("H5").Value = 741.8570999999999
$ws.Range("I5").Value = 477
$ws.Range("J5").Value = 786
$ws.Range("K5").Value = 1431
$ws.Range("L5").Value = 2358
$ws.Range("M5").Value = -1319
$ws.Range("N5").Value = -2582

$ws.Range("H122").Value = 2020
$ws.Range("I122").Value = 2028.8572
$ws.Range("J122").Value = 1999.3334
$ws.Range("K122").Value = 18259.7148
$ws.Range("L122").Value = 17994.0006
$ws.Range("M122").Value = -15809.7148
$ws.Range("N122").Value = -22894.0006

$ws.Range("H135").Value = 741.8570999999999
$ws.Range("I135").Value = 477
$ws.Range("J135").Value = 786
$ws.Range("K135").Value = 4293
$ws.Range("L135").Value = 7074
$ws.Range("M135").Value = -1758
$ws.Range("N135").Value = -12144

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 21888.5
$ws.Range("I59").Value = 2000
$ws.Range("J59").Value = 28518
$ws.Range("K59").Value = 2000
$ws.Range("L59").Value = 28518
$ws.Range("M59").Value = -1417
$ws.Range("N59").Value = -29684

$ws.Range("H114").Value = 50000
$ws.Range("J114").Value = 50000
$ws.Range("L114").Value = 50000
$ws.Range("N114").Value = -58678

$ws.Range("H132").Value = 93709.91
$ws.Range("I132").Value = 202270.8
$ws.Range("J132").Value = 3242.5
$ws.Range("K132").Value = 606812.3999999999
$ws.Range("L132").Value = 9727.5
$ws.Range("M132").Value = -604282.3999999999
$ws.Range("N132").Value = -14787.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3003.0386
$ws.Range("I132").Value = 3282.353
$ws.Range("J132").Value = 2475.4443
$ws.Range("K132").Value = 9847.059000000001
$ws.Range("L132").Value = 7426.3329
$ws.Range("M132").Value = -7317.059000000001
$ws.Range("N132").Value = -12486.3329

$ws.Range("H136").Value = 6207.2
$ws.Range("I136").Value = 10801.6
$ws.Range("J136").Value = 1612.8
$ws.Range("K136").Value = 32404.8
$ws.Range("L136").Value = 4838.4
$ws.Range("M136").Value = -29854.8
$ws.Range("N136").Value = -9938.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 101464.5
$ws.Range("J46").Value = 101464.5
$ws.Range("L46").Value = 101464.5
$ws.Range("N46").Value = -101926.5

$ws.Range("H132").Value = 1541.4
$ws.Range("I132").Value = 1239.8572
$ws.Range("J132").Value = 3124.5
$ws.Range("K132").Value = 3719.5716
$ws.Range("L132").Value = 9373.5
$ws.Range("M132").Value = -1189.5716
$ws.Range("N132").Value = -14433.5

$ws.Range("H134").Value = 101464.5
$ws.Range("J134").Value = 101464.5
$ws.Range("L134").Value = 304393.5
$ws.Range("N134").Value = -309463.5

$ws.Range("H136").Value = 5143.1934
$ws.Range("I136").Value = 6471.2607
$ws.Range("J136").Value = 1325
$ws.Range("K136").Value = 19413.7821
$ws.Range("L136").Value = 3975
$ws.Range("M136").Value = -16863.7821
$ws.Range("N136").Value = -9075
